# Add two new IN/OUT/HOURS:MINS time-tracking entries (rows 40-43) to the
# bottom of the time log, following the same repeating pattern used
# throughout the sheet: a bold/bordered header row (IN / OUT / HOURS:MINS)
# followed by a data row (date, time in, time out, duration).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (and, for column A, the already-textual date value)
# from the last existing header+data pair (rows 38:39) onto the two new
# pairs, so the bold font / thin-border styling used by every other entry
# in the sheet is preserved, and so the date is carried over as literal
# text instead of being re-parsed into a date serial number.
$ws.Range("B38:D39").Copy($ws.Range("B40"))
$ws.Range("A39").Copy($ws.Range("A41"))

$ws.Range("B38:D39").Copy($ws.Range("B42"))
$ws.Range("A39").Copy($ws.Range("A43"))

$excel.CutCopyMode = 0

# Row 40 - header
$ws.Range("B40").Value2 = "IN"
$ws.Range("C40").Value2 = "OUT"
$ws.Range("D40").Value2 = "HOURS:MINS"

# Row 41 - data (A41 already holds the correct text "2020-10-22" from the copy above)
$ws.Range("B41").Value2 = "15:26:31"
$ws.Range("C41").Value2 = "15:26:32"
$ws.Range("D41").Value2 = "0:0"

# Row 42 - header
$ws.Range("B42").Value2 = "IN"
$ws.Range("C42").Value2 = "OUT"
$ws.Range("D42").Value2 = "HOURS:MINS"

# Row 43 - data (A43 already holds the correct text "2020-10-22" from the copy above)
$ws.Range("B43").Value2 = "15:26:31"
$ws.Range("C43").Value2 = "15:26:33"
$ws.Range("D43").Value2 = "0:0"

Write-Host "New used range: $($ws.UsedRange.Address())"
